# KIBON-1621: Erlaeuterungen und Abwesenheit Info Text
$d = $word.ActiveDocument

# 1. Remove the _GoBack bookmark (bookmarkStart/bookmarkEnd pair) in the first paragraph.
try {
    $gb = $d.Bookmarks.Item("_GoBack")
    $gb.Delete()
} catch {
}

# 2. Strip the " (en pour cent)" suffix from the three bold heading paragraphs.
#    ("Taux de prise en charge effectif", "... accordé", "... subventionné")
$headings = @(
    "Taux de prise en charge effectif (en pour cent)",
    "Taux de prise en charge accordé (en pour cent)",
    "Taux de prise en charge subventionné (en pour cent)"
)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13)
    foreach ($h in $headings) {
        if ($t -eq $h) {
            [void]$p.Range.Find.Execute(" (en pour cent)", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)
        }
    }
}

# 3. Insert the new explanatory paragraph after "Il s'agit du taux de prise en charge convenu avec l'institution."
#    (right before "Taux de prise en charge accordé")
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13) -eq "Il s’agit du taux de prise en charge convenu avec l’institution.  ") {
        $p.Range.InsertParagraphAfter()
        $newP = $d.Paragraphs.Item($i + 1)
        $newP.Range.Text = "Les taux sont indiqués en pour cent (garderies) ou en heures (familles d’accueil). Une durée de prise en charge de 220 heures par mois chez des parents de jour correspond à un taux de 100%. La réduction du taux de prise en charge implique une diminution linéaire de la durée de prise en charge."
        break
    }
}

# 4. Insert the new paragraph about registering effective rate adaptations after
#    "Ce taux correspond au taux de prise en charge effectif ... seul ce dernier est subventionné."
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13) -eq "Ce taux correspond au taux de prise en charge effectif jusqu’à concurrence du maximum. Si le taux effectif est plus élevé que le taux accordé, seul ce dernier est subventionné.") {
        $p.Range.InsertParagraphAfter()
        $newP = $d.Paragraphs.Item($i + 1)
        $newP.Range.Text = "La structure d’accueil a jusqu’à la fin de la période tarifaire pour enregistrer les adaptations du taux de prise en charge effectif dans le cadre du taux de prise en charge accordé."
        break
    }
}

# 5. Expand the "Il s'agit de la réduction..." paragraph with the additional details.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13) -eq "Il s’agit de la réduction accordée, calculée sur la base du taux de prise en charge subventionné et du revenu déterminant. ") {
        $r = $p.Range
        [void]$r.MoveEnd(1, -1)
        $r.Text = "Il s’agit de la réduction accordée, calculée sur la base de la subvention maximale par unité de prise en charge, du taux de prise en charge subventionné, du revenu déterminant et d’un éventuel forfait pour frais de garde extraordinaires.  "
        break
    }
}
